$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newDate = 43199.562142037

$ws.Cells.Item(3, 2).Value = $newDate
$ws.Cells.Item(3, 3).Value = $newDate
$ws.Cells.Item(3, 4).Value = $newDate
$ws.Cells.Item(3, 5).Value = $newDate
$ws.Cells.Item(3, 6).Value = $newDate
$ws.Cells.Item(3, 7).Value = $newDate
